$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.196.92"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").Value = "1.914.35"
$ws.Range("E3").Value = "  +1.76%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "327.94"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("D6").Value = "1.006"

$ws.Range("D7").Value = "0.4665"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").Value = "0.3969"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.08028"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.014"
$ws.Range("E10").Value = "  +3.29%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "22.23"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.919.35"
$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "7.159"
$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.814"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.06967"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "89.39"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001021"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "17.35"
$ws.Range("E19").Value = "  +2.19%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "29.215.86"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.413"
$ws.Range("E22").Value = "  +1.16%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "11.17"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.171.22"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.058"
$ws.Range("E25").Value = "  -3.30%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "155.66"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.71"
$ws.Range("E27").Value = "  +1.53%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "5.942"
$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.039"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "120.99"
$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.09409"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.9459"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.377"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.354"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "3.280"
$ws.Range("E35").Value = "  -2.08%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.05894"
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "8.178"
$ws.Range("E37").Value = "  +3.51%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.175"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02111"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5859"
$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.1822"
$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "10.08"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.325"
$ws.Range("E44").Value = "  +12.04%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5477"
$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("D46").Value = "11.99"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.07213"
$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.927"
$ws.Range("E48").Value = "  +4.34%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "1.129"
$ws.Range("E49").Value = "  -4.30%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "113.28"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").Value = "2.501"
$ws.Range("E51").Value = "  +5.33%  "
